$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S6").Value = 1664
$ws.Range("B8").Value = "'701"
$ws.Range("B8").Style = "Normal"
$ws.Range("B9").Value = "'704"
$ws.Range("B9").Style = "Normal"
$ws.Range("B10").Value = "'707"
$ws.Range("B10").Style = "Normal"
$ws.Range("B11").Value = "'710"
$ws.Range("B11").Style = "Normal"
$ws.Range("B12").Value = "'457"
$ws.Range("B12").Style = "Normal"
$ws.Range("B13").Value = "'713"
$ws.Range("B13").Style = "Normal"
$ws.Range("B14").Value = "'716"
$ws.Range("B14").Style = "Normal"
$ws.Range("B15").Value = "'719"
$ws.Range("B15").Style = "Normal"
$ws.Range("B16").Value = "'467"
$ws.Range("B16").Style = "Normal"
$ws.Range("B17").Value = "'722"
$ws.Range("B17").Style = "Normal"
$ws.Range("B18").Value = "'725"
$ws.Range("B18").Style = "Normal"
$ws.Range("B19").Value = "'474"
$ws.Range("B19").Style = "Normal"
$ws.Range("B20").Value = "'728"
$ws.Range("B20").Style = "Normal"
$ws.Range("B21").Value = "'731"
$ws.Range("B21").Style = "Normal"
$ws.Range("B22").Value = "'734"
$ws.Range("B22").Style = "Normal"
$ws.Range("B23").Value = "'484"
$ws.Range("B23").Style = "Normal"
$ws.Range("B24").Value = "'737"
$ws.Range("B24").Style = "Normal"
$ws.Range("B25").Value = "'740"
$ws.Range("B25").Style = "Normal"
$ws.Range("B26").Value = "'491"
$ws.Range("B26").Style = "Normal"
$ws.Range("B27").Value = "'743"
$ws.Range("B27").Style = "Normal"
$ws.Range("B28").Value = "'746"
$ws.Range("B28").Style = "Normal"
$ws.Range("B29").Value = "'498"
$ws.Range("B29").Style = "Normal"
$ws.Range("B30").Value = "'749"
$ws.Range("B30").Style = "Normal"
$ws.Range("B31").Value = "'752"
$ws.Range("B31").Style = "Normal"
$ws.Range("B32").Value = "'505"
$ws.Range("B32").Style = "Normal"
$ws.Range("B33").Value = "'755"
$ws.Range("B33").Style = "Normal"
$ws.Range("B34").Value = "'758"
$ws.Range("B34").Style = "Normal"
$ws.Range("B35").Value = "'512"
$ws.Range("B35").Style = "Normal"
$ws.Range("B36").Value = "'761"
$ws.Range("B36").Style = "Normal"
$ws.Range("B37").Value = "'764"
$ws.Range("B37").Style = "Normal"
$ws.Range("B38").Value = "'519"
$ws.Range("B38").Style = "Normal"
$ws.Range("B39").Value = "'767"
$ws.Range("B39").Style = "Normal"
$ws.Range("B40").Value = "'770"
$ws.Range("B40").Style = "Normal"
$ws.Range("B41").Value = "'526"
$ws.Range("B41").Style = "Normal"
$ws.Range("B42").Value = "'773"
$ws.Range("B42").Style = "Normal"
$ws.Range("B43").Value = "'776"
$ws.Range("B43").Style = "Normal"
$ws.Range("B44").Value = "'533"
$ws.Range("B44").Style = "Normal"
$ws.Range("B45").Value = "'779"
$ws.Range("B45").Style = "Normal"
$ws.Range("B46").Value = "'782"
$ws.Range("B46").Style = "Normal"
$ws.Range("B47").Value = "'785"
$ws.Range("B47").Style = "Normal"
$ws.Range("B48").Value = "'543"
$ws.Range("B48").Style = "Normal"
$ws.Range("B49").Value = "'788"
$ws.Range("B49").Style = "Normal"
$ws.Range("B50").Value = "'791"
$ws.Range("B50").Style = "Normal"
$ws.Range("B51").Value = "'550"
$ws.Range("B51").Style = "Normal"
$ws.Range("B52").Value = "'794"
$ws.Range("B52").Style = "Normal"
$ws.Range("B53").Value = "'797"
$ws.Range("B53").Style = "Normal"
$ws.Range("B54").Value = "'557"
$ws.Range("B54").Style = "Normal"
$ws.Range("B55").Value = "'800"
$ws.Range("B55").Style = "Normal"
$ws.Range("B56").Value = "'803"
$ws.Range("B56").Style = "Normal"
$ws.Range("B57").Value = "'564"
$ws.Range("B57").Style = "Normal"
$ws.Range("B58").Value = "'806"
$ws.Range("B58").Style = "Normal"
$ws.Range("B59").Value = "'809"
$ws.Range("B59").Style = "Normal"
$ws.Range("B60").Value = "'571"
$ws.Range("B60").Style = "Normal"
$ws.Range("B61").Value = "'812"
$ws.Range("B61").Style = "Normal"
$ws.Range("B62").Value = "'815"
$ws.Range("B62").Style = "Normal"
$ws.Range("B63").Value = "'578"
$ws.Range("B63").Style = "Normal"
$ws.Range("B64").Value = "'818"
$ws.Range("B64").Style = "Normal"
$ws.Range("B65").Value = "'821"
$ws.Range("B65").Style = "Normal"
$ws.Range("B66").Value = "'585"
$ws.Range("B66").Style = "Normal"
$ws.Range("B67").Value = "'824"
$ws.Range("B67").Style = "Normal"
$ws.Range("B68").Value = "'827"
$ws.Range("B68").Style = "Normal"
$ws.Range("B69").Value = "'592"
$ws.Range("B69").Style = "Normal"
$ws.Range("B70").Value = "'830"
$ws.Range("B70").Style = "Normal"
$ws.Range("B71").Value = "'833"
$ws.Range("B71").Style = "Normal"
$ws.Range("B72").Value = "'836"
$ws.Range("B72").Style = "Normal"
$ws.Range("B73").Value = "'602"
$ws.Range("B73").Style = "Normal"
$ws.Range("B74").Value = "'839"
$ws.Range("B74").Style = "Normal"
$ws.Range("B75").Value = "'842"
$ws.Range("B75").Style = "Normal"
$ws.Range("B76").Value = "'609"
$ws.Range("B76").Style = "Normal"
$ws.Range("B77").Value = "'845"
$ws.Range("B77").Style = "Normal"
$ws.Range("B78").Value = "'848"
$ws.Range("B78").Style = "Normal"
$ws.Range("B79").Value = "'616"
$ws.Range("B79").Style = "Normal"
$ws.Range("B80").Value = "'851"
$ws.Range("B80").Style = "Normal"
$ws.Range("B81").Value = "'854"
$ws.Range("B81").Style = "Normal"
$ws.Range("B82").Value = "'623"
$ws.Range("B82").Style = "Normal"
$ws.Range("B83").Value = "'857"
$ws.Range("B83").Style = "Normal"
$ws.Range("B84").Value = "'860"
$ws.Range("B84").Style = "Normal"
$ws.Range("B85").Value = "'863"
$ws.Range("B85").Style = "Normal"
$ws.Range("B86").Value = "'866"
$ws.Range("B86").Style = "Normal"
$ws.Range("B87").Value = "'869"
$ws.Range("B87").Style = "Normal"
$ws.Range("B88").Value = "'872"
$ws.Range("B88").Style = "Normal"
$ws.Range("B90").Value = "'875"
$ws.Range("B90").Style = "Normal"
$ws.Range("B91").Value = "'651"
$ws.Range("B91").Style = "Normal"
$ws.Range("B92").Value = "'878"
$ws.Range("B92").Style = "Normal"
$ws.Range("B93").Value = "'881"
$ws.Range("B93").Style = "Normal"
$ws.Range("B94").Value = "'884"
$ws.Range("B94").Style = "Normal"
$ws.Range("B95").Value = "'887"
$ws.Range("B95").Style = "Normal"
$ws.Range("B96").Value = "'890"
$ws.Range("B96").Style = "Normal"
$ws.Range("B97").Value = "'893"
$ws.Range("B97").Style = "Normal"
$ws.Range("B99").Value = "'896"
$ws.Range("B99").Style = "Normal"
$ws.Range("B100").Value = "'680"
$ws.Range("B100").Style = "Normal"
$ws.Range("B101").Value = "'899"
$ws.Range("B101").Style = "Normal"
$ws.Range("B102").Value = "'902"
$ws.Range("B102").Style = "Normal"
$ws.Range("B103").Value = "'905"
$ws.Range("B103").Style = "Normal"
$ws.Range("B104").Value = "'908"
$ws.Range("B104").Style = "Normal"
$ws.Range("B105").Value = "'911"
$ws.Range("B105").Style = "Normal"
$ws.Range("B106").Value = "'914"
$ws.Range("B106").Style = "Normal"
$ws.Range("S114").Value = 186810
